$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.118515014648438
$ws.Range("B1").Value = 2.814093112945557
$ws.Range("C1").Value = 6.866245269775391
$ws.Range("D1").Value = 2.022721529006958
$ws.Range("E1").Value = 1.079953670501709
